$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting used by the other data rows (green fill, left aligned)
# onto the new rows that will receive data (rows 29-35).
$ws.Range("A2:C2").Copy()
$ws.Range("A29:C35").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Developer / File Name / VBA File Strings Length for the new entries
$ws.Range("A29").Value = "Vahid"
$ws.Range("B29").Value = "TubeDefDisc"
$ws.Range("C29").Value = 133

$ws.Range("A30").Value = "Vahid"
$ws.Range("B30").Value = "cClashData"
$ws.Range("C30").Value = 52

$ws.Range("A31").Value = "Vahid"
$ws.Range("B31").Value = "cLogEntry"
$ws.Range("C31").Value = 39

$ws.Range("A32").Value = "Vahid"
$ws.Range("B32").Value = "cPlotSupt"
$ws.Range("C32").Value = 275

$ws.Range("A33").Value = "Vahid"
$ws.Range("B33").Value = "cSteelDisc"
$ws.Range("C33").Value = 196

$ws.Range("A34").Value = "Vahid"
$ws.Range("B34").Value = "cSuptPoints"
$ws.Range("C34").Value = 399

$ws.Range("A35").Value = "Vahid"
$ws.Range("B35").Value = "cTubeDef"
$ws.Range("C35").Value = 174

# Update view: move the active selection to D19
$ws.Range("D19").Select()
